# Add a new model row (Gemma-7B-Instruct) as row 8 to the results sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Gemma-7B-Instruct"
$ws.Range("B8").Value = "0.53 ± 0.5"
$ws.Range("C8").Value = "-0.13 ± 0.89"
$ws.Range("D8").Value = "0.22 ± 0.53"
$ws.Range("E8").Value = "0.01 ± 0.01"
$ws.Range("F8").Value = "0.09 ± 0.1"
$ws.Range("G8").Value = "0.01 ± 0.03"
$ws.Range("H8").Value = "0.08 ± 0.09"
$ws.Range("I8").Value = "0.09 ± 0.1"
$ws.Range("J8").Value = "0.48 ± 0.42"
$ws.Range("K8").Value = "0.49 ± 0.42"
$ws.Range("L8").Value = "0.49 ± 0.42"
$ws.Range("M8").Value = "0.48 ± 0.42"
$ws.Range("N8").Value = "0.56 ± 0.49"
$ws.Range("O8").Value = "0.08 ± 0.1"
$ws.Range("P8").Value = "0.39 ± 0.34"
$ws.Range("Q8").Value = "5.17 ± 1.03"
$ws.Range("R8").Value = "0.425 ± 0.00"
$ws.Range("S8").Value = "0.54 ± 0.46"
$ws.Range("T8").Value = "0.57 ± 0.49"
$ws.Range("U8").Value = "2.1 ± 1.97"
$ws.Range("V8").Value = "0.55 ± 0.49"
$ws.Range("W8").Value = "0.54 ± 0.46"
$ws.Range("X8").Value = "0.86 ± 0.76"
